# Do not exempt process emissions from carbon tax by default.
# The control-lever toggle on the "BEPEfCT" sheet (cell B2) is switched
# from 1 (exempt) to 0 (not exempt / taxed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")
$ws.Range("B2").Value = 0
